$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value would otherwise be auto-parsed as a number by Excel;
# force them to Text format first so they stay text, matching the source data (inlineStr).
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Apply the updated values (price / volume / coin name / link changes).
$ws.Range("D2").Value = "23.326.81"
$ws.Range("D3").Value = "1.624.28"
$ws.Range("E3").Value = "  -0.93%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("D6").Value = "302.64"
$ws.Range("E6").Value = "  -0.81%  "
$ws.Range("D7").Value = "0.3749"
$ws.Range("E7").Value = "  +0.57%  "
$ws.Range("D8").Value = "0.3621"
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "51.27"
$ws.Range("E9").Value = "  -1.63%  "
$ws.Range("E10").Value = "  +0.24%  "
$ws.Range("D11").Value = "1.223"
$ws.Range("E11").Value = "  -2.26%  "
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("E13").Value = "  -2.21%  "
$ws.Range("D14").Value = "6.471"
$ws.Range("E14").Value = "  -1.81%  "
$ws.Range("D15").Value = "0.00001236"
$ws.Range("E15").Value = "  -2.49%  "
$ws.Range("D16").Value = "7.266"
$ws.Range("E16").Value = "  -0.22%  "
$ws.Range("D17").Value = "1.618.52"
$ws.Range("E17").Value = "  -0.91%  "
$ws.Range("D18").Value = "93.94"
$ws.Range("E18").Value = "  -0.39%  "
$ws.Range("D19").Value = "0.06927"
$ws.Range("E19").Value = "  +0.68%  "
$ws.Range("D20").Value = "17.53"
$ws.Range("E20").Value = "  -3.32%  "
$ws.Range("D21").Value = "6.526"
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("E23").Value = "  -1.63%  "
$ws.Range("D24").Value = "23.328.78"
$ws.Range("E24").Value = "  -0.41%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "2.468"
$ws.Range("E25").Value = "  +2.64%  "
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").Value = "3.082"
$ws.Range("E26").Value = "  +1.54%  "
$ws.Range("D27").Value = "21.11"
$ws.Range("E27").Value = "  -0.51%  "
$ws.Range("D28").Value = "150.60"
$ws.Range("E28").Value = "  -0.77%  "
$ws.Range("D29").Value = "5.268"
$ws.Range("E29").Value = "  -0.97%  "
$ws.Range("E30").Value = "  -2.20%  "
$ws.Range("D31").Value = "1.798.72"
$ws.Range("E31").Value = "  -0.83%  "
$ws.Range("D32").Value = "6.738"
$ws.Range("E32").Value = "  -0.14%  "
$ws.Range("D33").Value = "2.167"
$ws.Range("E33").Value = "  -5.40%  "
$ws.Range("D34").Value = "1.062"
$ws.Range("E34").Value = "  +11.40%  "
$ws.Range("D35").Value = "11.22"
$ws.Range("E35").Value = "  +8.72%  "
$ws.Range("D36").Value = "0.02751"
$ws.Range("E36").Value = "  -3.22%  "
$ws.Range("D37").Value = "0.08754"
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("E38").Value = "  -1.49%  "
$ws.Range("D39").Value = "0.07101"
$ws.Range("E39").Value = "  -1.77%  "
$ws.Range("D40").Value = "5.992"
$ws.Range("E40").Value = "  -0.97%  "
$ws.Range("D41").Value = "0.6969"
$ws.Range("E41").Value = "  -1.07%  "
$ws.Range("D42").Value = "1.334"
$ws.Range("E42").Value = "  -2.76%  "
$ws.Range("D43").Value = "16.07"
$ws.Range("E43").Value = "  +0.49%  "
$ws.Range("D44").Value = "12.05"
$ws.Range("E44").Value = "  -3.16%  "
$ws.Range("D45").Value = "0.6461"
$ws.Range("E45").Value = "  -0.70%  "
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").Value = "3.957"
$ws.Range("E47").Value = "  -1.22%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "2.262"
$ws.Range("E48").Value = "  -2.73%  "
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("D50").Value = "125.67"
$ws.Range("E50").Value = "  -1.93%  "
$ws.Range("E51").Value = "  -1.38%  "
